$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STE")

# Row 4 - Inventory
$ws.Range("B4").Value = 294000000.0
$ws.Range("C4").Value = 279000000.0
$ws.Range("D4").Value = 277000000.0
$ws.Range("E4").Value = 248000000.0
$ws.Range("F4").Value = 252000000.0

# Row 12 - Accounts Payable
$ws.Range("B12").Value = 134000000.0
$ws.Range("C12").Value = 125000000.0
$ws.Range("D12").Value = 134000000.0
$ws.Range("E12").Value = 149000000.0
$ws.Range("F12").Value = 140000000.0

# Row 19 - Long Term Tax Liability (Deferred)
$ws.Range("B19").Value = 264000000.0
$ws.Range("C19").Value = 165000000.0
$ws.Range("D19").Value = 161000000.0
$ws.Range("E19").Value = 161000000.0
$ws.Range("F19").Value = 156000000.0
